$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted cells (Price/Volume columns) keep their exact string
# representation (e.g. trailing zeros, multi-dot numbers) instead of being
# auto-converted to numbers by Excel when the value is assigned.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.346.87'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +4.64%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.717.27'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.79%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9979'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.35'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.51%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2639'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06233'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.91%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.712.12'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07084'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.27'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +5.43%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5931'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.427'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.20'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9991'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9986'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.338.58'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +4.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006816'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.43%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.925.80'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.88%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.566'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +5.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.851'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +3.99%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.355'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '135.69'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.20'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.404'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +2.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.768'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +7.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.53'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.037'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.691'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +3.68%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07772'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.05%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.609'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6234'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9720'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.63%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9167'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +7.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '112.10'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +13.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.422'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -6.32%  '
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.000'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.10%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.905'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +6.05%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.3818'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.16%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.163'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +11.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1145'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +4.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.257'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05292'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '30.72'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +5.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.667'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +6.38%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.48%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3387'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.73%  '
